# Edit script for B0C5BQ1HM6_po_data.xlsx
# 1. Rename header cell B1 on "Weekly Quantity" to "Weekly_PO_Qty"
# 2. Rename header cell B1 on "Monthly Trend" to "Monthly_PO_Qty"
# 3. Add a new "PO Forecast" worksheet with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Copy header formatting (bold, centered, bordered) from "Weekly Quantity" row 1
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-formatted style from "Weekly Quantity" column A onto the new sheet's ds column
$wsWeekly.Range("A2:A38").Copy()
$wsForecast.Range("A2:A46").PasteSpecial(-4122)

# Fill in the forecast data (rows 2-46)
$arr = New-Object 'object[,]' 45,4
$arr[0,0] = 45067.99999999999; $arr[0,1] = 377; $arr[0,2] = -167.6201035971666; $arr[0,3] = 895.9343489773556
$arr[1,0] = 45088.99999999999; $arr[1,1] = 375; $arr[1,2] = -149.5526682201012; $arr[1,3] = 934.1962061527547
$arr[2,0] = 45095.99999999999; $arr[2,1] = 375; $arr[2,2] = -160.5059425075265; $arr[2,3] = 855.1558878676997
$arr[3,0] = 45102.99999999999; $arr[3,1] = 374; $arr[3,2] = -164.210164477648; $arr[3,3] = 923.1209601553389
$arr[4,0] = 45109.99999999999; $arr[4,1] = 374; $arr[4,2] = -161.7013735420639; $arr[4,3] = 912.2807578878048
$arr[5,0] = 45116.99999999999; $arr[5,1] = 373; $arr[5,2] = -196.408313147491; $arr[5,3] = 936.9213129534812
$arr[6,0] = 45123.99999999999; $arr[6,1] = 373; $arr[6,2] = -201.2726136946428; $arr[6,3] = 933.2241788391301
$arr[7,0] = 45130.99999999999; $arr[7,1] = 372; $arr[7,2] = -180.1303586010583; $arr[7,3] = 943.9439542685134
$arr[8,0] = 45137.99999999999; $arr[8,1] = 372; $arr[8,2] = -173.9699601615456; $arr[8,3] = 900.6778233229171
$arr[9,0] = 45151.99999999999; $arr[9,1] = 371; $arr[9,2] = -171.7216516474463; $arr[9,3] = 889.720446561116
$arr[10,0] = 45158.99999999999; $arr[10,1] = 370; $arr[10,2] = -169.7470307266962; $arr[10,3] = 883.2509509398243
$arr[11,0] = 45172.99999999999; $arr[11,1] = 369; $arr[11,2] = -168.8468315567614; $arr[11,3] = 924.0081826539438
$arr[12,0] = 45186.99999999999; $arr[12,1] = 368; $arr[12,2] = -136.4426745778821; $arr[12,3] = 957.8230232608395
$arr[13,0] = 45193.99999999999; $arr[13,1] = 368; $arr[13,2] = -171.4405645503573; $arr[13,3] = 896.4963622374023
$arr[14,0] = 45200.99999999999; $arr[14,1] = 367; $arr[14,2] = -168.2168377777299; $arr[14,3] = 891.5761228312274
$arr[15,0] = 45207.99999999999; $arr[15,1] = 367; $arr[15,2] = -200.6822293655187; $arr[15,3] = 891.4496829729277
$arr[16,0] = 45221.99999999999; $arr[16,1] = 366; $arr[16,2] = -174.8038718268915; $arr[16,3] = 935.8370435624231
$arr[17,0] = 45235.99999999999; $arr[17,1] = 365; $arr[17,2] = -206.3606404848967; $arr[17,3] = 924.3849478300446
$arr[18,0] = 45270.99999999999; $arr[18,1] = 362; $arr[18,2] = -161.5705388499583; $arr[18,3] = 916.7563720458726
$arr[19,0] = 45277.99999999999; $arr[19,1] = 362; $arr[19,2] = -204.1019130478203; $arr[19,3] = 952.0788718988316
$arr[20,0] = 45298.99999999999; $arr[20,1] = 360; $arr[20,2] = -179.7199153645465; $arr[20,3] = 916.5531771526192
$arr[21,0] = 45305.99999999999; $arr[21,1] = 360; $arr[21,2] = -174.6933223349998; $arr[21,3] = 919.1107694904301
$arr[22,0] = 45312.99999999999; $arr[22,1] = 359; $arr[22,2] = -217.8348304726932; $arr[22,3] = 874.7336737211343
$arr[23,0] = 45326.99999999999; $arr[23,1] = 358; $arr[23,2] = -185.6487056668201; $arr[23,3] = 856.3476966544197
$arr[24,0] = 45333.99999999999; $arr[24,1] = 358; $arr[24,2] = -187.7737504174749; $arr[24,3] = 939.4779741038634
$arr[25,0] = 45347.99999999999; $arr[25,1] = 357; $arr[25,2] = -203.9972742803984; $arr[25,3] = 869.0449370182019
$arr[26,0] = 45354.99999999999; $arr[26,1] = 356; $arr[26,2] = -175.5965650698359; $arr[26,3] = 891.4855508382385
$arr[27,0] = 45361.99999999999; $arr[27,1] = 356; $arr[27,2] = -166.6734261414533; $arr[27,3] = 912.1060039759194
$arr[28,0] = 45375.99999999999; $arr[28,1] = 355; $arr[28,2] = -197.2841786426922; $arr[28,3] = 904.8897696798249
$arr[29,0] = 45389.99999999999; $arr[29,1] = 354; $arr[29,2] = -190.2060748504325; $arr[29,3] = 926.2544792042419
$arr[30,0] = 45396.99999999999; $arr[30,1] = 353; $arr[30,2] = -176.1675280515049; $arr[30,3] = 884.669526874223
$arr[31,0] = 45480.99999999999; $arr[31,1] = 347; $arr[31,2] = -209.0262254332114; $arr[31,3] = 864.4999302174602
$arr[32,0] = 45536.99999999999; $arr[32,1] = 343; $arr[32,2] = -209.3306155616106; $arr[32,3] = 867.0496288862876
$arr[33,0] = 45543.99999999999; $arr[33,1] = 342; $arr[33,2] = -220.782731865289; $arr[33,3] = 849.8316142590845
$arr[34,0] = 45550.99999999999; $arr[34,1] = 342; $arr[34,2] = -210.494817077649; $arr[34,3] = 883.6147332750969
$arr[35,0] = 45571.99999999999; $arr[35,1] = 340; $arr[35,2] = -225.6084687511702; $arr[35,3] = 877.2533594087135
$arr[36,0] = 45578.99999999999; $arr[36,1] = 340; $arr[36,2] = -208.5457721096875; $arr[36,3] = 877.0367071108279
$arr[37,0] = 45585.99999999999; $arr[37,1] = 339; $arr[37,2] = -215.4309665619556; $arr[37,3] = 848.7012859855914
$arr[38,0] = 45592.99999999999; $arr[38,1] = 339; $arr[38,2] = -189.2917712741222; $arr[38,3] = 868.2664667394394
$arr[39,0] = 45599.99999999999; $arr[39,1] = 338; $arr[39,2] = -204.7323535204872; $arr[39,3] = 885.2642605182391
$arr[40,0] = 45606.99999999999; $arr[40,1] = 338; $arr[40,2] = -221.3471607107979; $arr[40,3] = 845.4629079169182
$arr[41,0] = 45613.99999999999; $arr[41,1] = 337; $arr[41,2] = -241.4183256017047; $arr[41,3] = 850.2026107598673
$arr[42,0] = 45620.99999999999; $arr[42,1] = 337; $arr[42,2] = -230.0619005703225; $arr[42,3] = 846.9695703196379
$arr[43,0] = 45627.99999999999; $arr[43,1] = 336; $arr[43,2] = -186.1917378252672; $arr[43,3] = 864.371829583937
$arr[44,0] = 45634.99999999999; $arr[44,1] = 336; $arr[44,2] = -200.1735716737297; $arr[44,3] = 875.0654797713668

$wsForecast.Range("A2:D46").Value = $arr

# Restore the originally active sheet/selection
$wsForecast.Range("A1").Select()
$wsWeekly.Activate()
$wsWeekly.Range("A1").Select()
